{"js": "// Apply the tracked-change style edits described in the diff:\n//   1. Move the `_GoBack` bookmark from the end of the \"Jego podw\u0142adni...\"\n//      paragraph to the (empty) paragraph near the top of the document\n//      (the last of the blank paragraphs right before the \"Wykonanie:\" block).\n//   2. Fix a typo (\"pzyjmowania\" -> \"przyjmowania\") and insert a new clause\n//      (\"lub przekaza\u0107 j\u0105 do rozporz\u0105dzenia kierownikaowi jednostki\")\n//      in the \"W punkcie przyjmowania dokument\u00f3w...\" paragraph.\n//   3. Fix a few small wording/spelling issues in the \"Elektroniczny obieg\n//      dokument\u00f3w niesie...\" paragraph:\n//        \"po przez\"  -> \"poprzez\"\n//        \"cech\u0105 system jest\" -> \"cech\u0105 systemu jest\"\n//        \"kontrol\u0119 nad terminowo\u015bci\" -> \"kontroli nad terminowo\u015bci\u0105\"\n\n// --- 1. Move the _GoBack bookmark -----------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// The target paragraph is the empty paragraph that sits two paragraphs\n// above the \"Wykonanie:\" heading (i.e. the last of the blank paragraphs\n// right after the document title, just before the blank spacer paragraph\n// that precedes \"Wykonanie:\").\nlet wykonanieIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Wykonanie:\") {\n    wykonanieIndex = i;\n    break;\n  }\n}\n\nconst targetIndex = wykonanieIndex >= 2 ? wykonanieIndex - 2 : 11;\nconst targetParagraph = paragraphs.items[targetIndex];\ntargetParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Update the \"W punkcie przyjmowania dokument\u00f3w...\" paragraph -------\nconst oldSentence =\n  \"jak listonosz. Nie ma jej wi\u0119c w punkcie pzyjmowania interesant\u00f3w.\";\nconst newSentence =\n  \"jak listonosz lub przekaza\u0107 j\u0105 do rozporz\u0105dzenia kierownikaowi jednostki.\" +\n  \" Nie ma jej wi\u0119c w punkcie przyjmowania interesant\u00f3w.\";\n\nconst listonoszResults = context.document.body.search(oldSentence, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nlistonoszResults.load(\"text\");\nawait context.sync();\n\nif (listonoszResults.items.length > 0) {\n  listonoszResults.items[0].insertText(newSentence, \"Replace\");\n  await context.sync();\n}\n\n// --- 3. Update the \"Elektroniczny obieg dokument\u00f3w niesie...\" paragraph ---\nconst replacements = [\n  [\"po przez automatyczn\u0105\", \"poprzez automatyczn\u0105\"],\n  [\"Kolejn\u0105 cech\u0105 system jest\", \"Kolejn\u0105 cech\u0105 systemu jest\"],\n  [\n    \"posiadanie przejrzystej kontrol\u0119 nad terminowo\u015bci za\u0142atwiania\",\n    \"posiadanie przejrzystej kontroli nad terminowo\u015bci\u0105 za\u0142atwiania\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the tracked-change style edits described in the diff:\n#   1. Move the `_GoBack` bookmark from the end of the \"Jego podw\u0142adni...\"\n#      paragraph to the (empty) paragraph near the top of the document\n#      (the last of the blank paragraphs right before the \"Wykonanie:\" block).\n#   2. Fix a typo (\"pzyjmowania\" -> \"przyjmowania\") and insert a new clause\n#      (\"lub przekaza\u0107 j\u0105 do rozporz\u0105dzenia kierownikaowi jednostki\")\n#      in the \"W punkcie przyjmowania dokument\u00f3w...\" paragraph.\n#   3. Fix a few small wording/spelling issues in the \"Elektroniczny obieg\n#      dokument\u00f3w niesie...\" paragraph:\n#        \"po przez\"  -> \"poprzez\"\n#        \"cech\u0105 system jest\" -> \"cech\u0105 systemu jest\"\n#        \"kontrol\u0119 nad terminowo\u015bci\" -> \"kontroli nad terminowo\u015bci\u0105\"\n\n$d = $word.ActiveDocument\n\n# --- 1. Move the _GoBack bookmark ------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# The target paragraph is the empty paragraph that sits two paragraphs above\n# the \"Wykonanie:\" heading (i.e. the last of the blank paragraphs right after\n# the document title, just before the blank spacer paragraph that precedes\n# \"Wykonanie:\") - the same spot Word leaves the _GoBack bookmark when the\n# cursor was last positioned there before saving.\n$wykonanieIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq \"Wykonanie:\") {\n        $wykonanieIndex = $i\n        break\n    }\n}\n\nif ($wykonanieIndex -ge 3) {\n    $targetIndex = $wykonanieIndex - 2\n} else {\n    $targetIndex = 12\n}\n\n$targetParagraph = $d.Paragraphs.Item($targetIndex)\n$d.Bookmarks.Add(\"_GoBack\", $targetParagraph.Range)\n\n# --- 2. Update the \"W punkcie przyjmowania dokument\u00f3w...\" paragraph -------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"jak listonosz. Nie ma jej wi\u0119c w punkcie pzyjmowania interesant\u00f3w.\"\n$find.Replacement.Text = \"jak listonosz lub przekaza\u0107 j\u0105 do rozporz\u0105dzenia kierownikaowi jednostki. Nie ma jej wi\u0119c w punkcie przyjmowania interesant\u00f3w.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# --- 3. Update the \"Elektroniczny obieg dokument\u00f3w niesie...\" paragraph ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"po przez automatyczn\u0105\"\n$find2.Replacement.Text = \"poprzez automatyczn\u0105\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Kolejn\u0105 cech\u0105 system jest\"\n$find3.Replacement.Text = \"Kolejn\u0105 cech\u0105 systemu jest\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Replacement.ClearFormatting()\n$find4.Text = \"posiadanie przejrzystej kontrol\u0119 nad terminowo\u015bci za\u0142atwiania\"\n$find4.Replacement.Text = \"posiadanie przejrzystej kontroli nad terminowo\u015bci\u0105 za\u0142atwiania\"\n$find4.Execute($find4.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find4.Replacement.Text, 2)\n"}
